$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: new level entry (date, version, Entwicklung, Anna Franziska, info text) ---
# Seed formatting by copying the same-shaped row above (A14:E14 uses styles
# s=4 (date), s=3 (plain top-aligned), s=3, s=3, s=5 (wrap top-aligned)),
# then overwrite with the new row's values.
$ws.Range("A14:E14").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A15").Value2 = 42993
$ws.Range("B15").Value2 = "DiscordiaAgency_Demo_2017_09_15.exe"
$ws.Range("C15").Value2 = "Entwicklung"
$ws.Range("D15").Value2 = "Anna Franziska"
$ws.Range("E15").Value2 = 'insg. 6 Level eingebaut; Wachen können stationär sein; "globales" Alarmsystem: Wachen entdecken Leichen & geben größeren Suchradius an alle Wachen weiter'

$ws.Rows.Item(15).RowHeight = 75

# --- Row 16: new "Spielen" entry (version, Spielen, Urban) ---
$ws.Range("B14:D14").Copy()
$ws.Range("B16:D16").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B16").Value2 = "DiscordiaAgency_Demo_2017_09_15.exe"
$ws.Range("C16").Value2 = "Spielen"
$ws.Range("D16").Value2 = "Urban"

# --- Update selection to H15, matching the authored view state ---
$ws.Range("H15").Select()
